# "Code Check In 10 Oct"
# - Row 24 (ImportSIMs): fill in the previously-blank Description with
#   "Testing Required" and flip Runmode from Y to N.
# - Add a new row 25 for the ImportMojio test case (Description
#   "Testing Required", Runmode "Y"), matching the formatting already
#   used by the rest of the data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate row 24's formatting into row 25 before writing any values so
# the new row picks up the same cell style (bordered, non-header) as the
# rest of the table instead of the sheet default.
$ws.Range("A24:C24").Copy($ws.Range("A25:C25"))

# Row 24: ImportSIMs now needs testing and its Runmode flips to N.
$ws.Range("B24").Value = "Testing Required"
$ws.Range("C24").Value = "N"

# Row 25: new ImportMojio test case.
$ws.Range("A25").Value = "ImportMojio"
$ws.Range("B25").Value = "Testing Required"
$ws.Range("C25").Value = "Y"

$ws.Range("A25").Select()
